$wb = $excel.ActiveWorkbook

# Sheet 1 ("展览") F-column updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 233
$ws1.Range("F3").Value = 1080
$ws1.Range("F5").Value = 403
$ws1.Range("F6").Value = 75
$ws1.Range("F7").Value = 539
$ws1.Range("F8").Value = 56
$ws1.Range("F9").Value = 6709
$ws1.Range("F10").Value = 144
$ws1.Range("F15").Value = 1072
$ws1.Range("F16").Value = 16046
$ws1.Range("F18").Value = 33
$ws1.Range("F19").Value = 324
$ws1.Range("F20").Value = 172
$ws1.Range("F21").Value = 114
$ws1.Range("F22").Value = 11273
$ws1.Range("F23").Value = 2
$ws1.Range("F24").Value = 853
$ws1.Range("F25").Value = 4434
$ws1.Range("F26").Value = 294
$ws1.Range("F27").Value = 385
$ws1.Range("F28").Value = 38
$ws1.Range("F29").Value = 24
$ws1.Range("F31").Value = 136

# Sheet 4 ("全部类型") F-column updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 233
$ws4.Range("F3").Value = 1080
$ws4.Range("F5").Value = 403
$ws4.Range("F6").Value = 75
$ws4.Range("F7").Value = 539
$ws4.Range("F9").Value = 56
$ws4.Range("F10").Value = 6709
$ws4.Range("F11").Value = 144
$ws4.Range("F17").Value = 1072
$ws4.Range("F18").Value = 16047
$ws4.Range("F20").Value = 33
$ws4.Range("F21").Value = 324
$ws4.Range("F22").Value = 172
$ws4.Range("F23").Value = 114
$ws4.Range("F26").Value = 11273
$ws4.Range("F27").Value = 2
$ws4.Range("F28").Value = 854
$ws4.Range("F29").Value = 4434
$ws4.Range("F30").Value = 294
$ws4.Range("F31").Value = 385
$ws4.Range("F32").Value = 38
$ws4.Range("F33").Value = 24
$ws4.Range("F35").Value = 136
